$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (C) date on all three rows to 46064 (2026-02-11)
$ws.Range("C2").Value = 46064
$ws.Range("C3").Value = 46064
$ws.Range("C4").Value = 46064

# Swap data between row 3 and row 4 (Beteckning, Datum, Area)
$ws.Range("A3").Value = "A 36713-2023"
$ws.Range("B3").Value = 45153
$ws.Range("G3").Value = 0.7

$ws.Range("A4").Value = "A 35536-2025"
$ws.Range("B4").Value = 45856
$ws.Range("G4").Value = 4
